# Update MRE device names in the "mre_costs" sheet: drop the possessive
# "'s" from each manufacturer name (e.g. "Nova Innovation's M100-D" ->
# "Nova Innovation M100-D").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mre_costs")

$ws.Range("A2").Value = "Nova Innovation M100-D"
$ws.Range("A3").Value = "Tidal Energy Ltd Deltastream"
$ws.Range("A4").Value = "Alstom DeepGen"
$ws.Range("A5").Value = "Orbital Marine Power SR200"
$ws.Range("A6").Value = "Orbital Marine Power O2"
$ws.Range("A7").Value = "Andritz Hydro HS1500"
$ws.Range("A8").Value = "SIMEC AR1500"

# Match the saved selection/active cell on this sheet.
$ws.Activate()
$ws.Range("A6").Select() | Out-Null
